$wb = $excel.ActiveWorkbook

# Each weekly timesheet runs Monday(4) Tuesday(5) Thursday(6) Friday(7)
# Saturday(8) Sunday(9) Total(11). A row for Wednesday was missing between
# Tuesday and Thursday - insert it and fix up the running-total formulas
# for the two rows that follow the insertion point (the engine shifts the
# remaining rows' relative formulas correctly on its own).
$weeks = @("Week 1","Week 2","Week 3","Week 4","Week 5","Week 6","Week 7","Week 8","Week 9","Week 10")

foreach ($wk in $weeks) {
    $ws = $wb.Worksheets.Item($wk)

    $ws.Rows.Item(6).Insert()

    $ws.Range("A6").Value = "Wednesday"
    $ws.Range("C6").Formula = "=B6+C5"
    $ws.Range("E6").Formula = "=D6+E5"
    $ws.Range("G6").Formula = "=F6+G5"
    $ws.Range("I6").Formula = "=H6+I5"

    $ws.Range("C7").Formula = "=B7+C6"
    $ws.Range("E7").Formula = "=D7+E6"
    $ws.Range("G7").Formula = "=F7+G6"
    $ws.Range("I7").Formula = "=H7+I6"

    $ws.Range("A6:I7").Select()
}

# Week 1 picked up a stray backtick left in G13 while the row was added.
$weekOne = $wb.Worksheets.Item("Week 1")
$weekOne.Range("G13").Value = "``"
$weekOne.Range("I7").Select()

# Week 2's selection ended up somewhere else after the edits were made.
$wb.Worksheets.Item("Week 2").Range("I7").Select()

# Totals sheet pulls the Sunday total from each weekly sheet - that row
# moved from 9 to 10 on every weekly sheet, so repoint the formulas.
$totals = $wb.Worksheets.Item("Totals")
$totals.Range("B4").Formula = "='Week 1'!C`$10"
$totals.Range("D4").Formula = "='Week 1'!E`$10"
$totals.Range("F4").Formula = "='Week 1'!G`$10"
$totals.Range("H4").Formula = "='Week 1'!I`$10"

$totals.Range("B5").Formula = "='Week 2'!C`$10"
$totals.Range("D5").Formula = "='Week 2'!E`$10"
$totals.Range("F5").Formula = "='Week 2'!G`$10"
$totals.Range("H5").Formula = "='Week 2'!I`$10"

$totals.Range("B6").Formula = "='Week 3'!C`$10"
$totals.Range("D6").Formula = "='Week 3'!E`$10"
$totals.Range("F6").Formula = "='Week 3'!G`$10"
$totals.Range("H6").Formula = "='Week 3'!I`$10"

$totals.Range("B7").Formula = "='Week 4'!C`$10"
$totals.Range("D7").Formula = "='Week 4'!E`$10"
$totals.Range("F7").Formula = "='Week 4'!G`$10"
$totals.Range("H7").Formula = "='Week 4'!I`$10"

$totals.Range("B8").Formula = "='Week 5'!C`$10"
$totals.Range("D8").Formula = "='Week 5'!E`$10"
$totals.Range("F8").Formula = "='Week 5'!G`$10"
$totals.Range("H8").Formula = "='Week 5'!I`$10"

$totals.Range("B9").Formula = "='Week 6'!C`$10"
$totals.Range("D9").Formula = "='Week 6'!E`$10"
$totals.Range("F9").Formula = "='Week 6'!G`$10"
$totals.Range("H9").Formula = "='Week 6'!I`$10"

$totals.Range("B10").Formula = "='Week 7'!C`$10"
$totals.Range("D10").Formula = "='Week 7'!E`$10"
$totals.Range("F10").Formula = "='Week 7'!G`$10"
$totals.Range("H10").Formula = "='Week 7'!I`$10"

$totals.Range("B11").Formula = "='Week 8'!C`$10"
$totals.Range("D11").Formula = "='Week 8'!E`$10"
$totals.Range("F11").Formula = "='Week 8'!G`$10"
$totals.Range("H11").Formula = "='Week 8'!I`$10"

$totals.Range("B12").Formula = "='Week 9'!C`$10"
$totals.Range("D12").Formula = "='Week 9'!E`$10"
$totals.Range("F12").Formula = "='Week 9'!G`$10"
$totals.Range("H12").Formula = "='Week 9'!I`$10"

$totals.Range("B13").Formula = "='Week 10'!C`$10"
$totals.Range("D13").Formula = "='Week 10'!E`$10"
$totals.Range("F13").Formula = "='Week 10'!G`$10"
$totals.Range("H13").Formula = "='Week 10'!I`$10"

$totals.Range("B6").Select()

# Leave the workbook with "Week 10" as the active tab/selection, matching
# where the author ended up after making all the edits.
$weekTen = $wb.Worksheets.Item("Week 10")
$weekTen.Activate()
$weekTen.Range("I14").Select()
